# Updated cryptos list (coinranking snapshot refresh).
#
# For every listed cell, write the refreshed value exactly as captured from
# the source data. The "Price" column (D) stores its values as plain text in
# this sheet (it keeps trailing zeros and thousand-dot groupings like
# "82.50" or "64.788.15"), so any replacement value in that column that
# LOOKS like an ordinary number is written with a leading apostrophe. That
# forces Excel to keep it as literal text instead of silently re-parsing it
# into a numeric cell, which would otherwise drop formatting (e.g.
# "82.50" -> 82.5, "12.00" -> 12). The quote-prefix styling that trick
# leaves behind is cleared right after via Style = "Normal" so the cell ends
# up as plain, unstyled text - same as every other text cell on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '65.093.35' },
    @{ Cell = 'E2'; Value = '  +0.12%  ' },
    @{ Cell = 'D3'; Value = '2.950.92' },
    @{ Cell = 'E3'; Value = '  -1.20%  ' },
    @{ Cell = 'E4'; Value = '  -0.17%  ' },
    @{ Cell = 'D5'; Value = '569.14' },
    @{ Cell = 'E5'; Value = '  -2.13%  ' },
    @{ Cell = 'D6'; Value = '159.82' },
    @{ Cell = 'E6'; Value = '  +3.97%  ' },
    @{ Cell = 'D7'; Value = '0.999' },
    @{ Cell = 'E7'; Value = '  -0.12%  ' },
    @{ Cell = 'D8'; Value = '0.519' },
    @{ Cell = 'E8'; Value = '  +0.81%  ' },
    @{ Cell = 'D9'; Value = '2.947.98' },
    @{ Cell = 'E9'; Value = '  -1.13%  ' },
    @{ Cell = 'D10'; Value = '6.68' },
    @{ Cell = 'E10'; Value = '  -4.54%  ' },
    @{ Cell = 'E11'; Value = '  +0.12%  ' },
    @{ Cell = 'D12'; Value = '0.456' },
    @{ Cell = 'E12'; Value = '  +2.20%  ' },
    @{ Cell = 'D13'; Value = '0.0000246' },
    @{ Cell = 'E13'; Value = '  +3.05%  ' },
    @{ Cell = 'D14'; Value = '34.10' },
    @{ Cell = 'E14'; Value = '  +0.32%  ' },
    @{ Cell = 'E15'; Value = '  -0.28%  ' },
    @{ Cell = 'D16'; Value = '65.195.66' },
    @{ Cell = 'E16'; Value = '  +0.24%  ' },
    @{ Cell = 'D17'; Value = '3.436.97' },
    @{ Cell = 'E17'; Value = '  -1.28%  ' },
    @{ Cell = 'D18'; Value = '6.93' },
    @{ Cell = 'E18'; Value = '  -0.01%  ' },
    @{ Cell = 'D19'; Value = '2.948.16' },
    @{ Cell = 'E19'; Value = '  -1.71%  ' },
    @{ Cell = 'D20'; Value = '446.05' },
    @{ Cell = 'E20'; Value = '  -0.50%  ' },
    @{ Cell = 'E21'; Value = '  +0.74%  ' },
    @{ Cell = 'D22'; Value = '0.681' },
    @{ Cell = 'E22'; Value = '  +0.10%  ' },
    @{ Cell = 'D23'; Value = '7.22' },
    @{ Cell = 'E23'; Value = '  -1.34%  ' },
    @{ Cell = 'D24'; Value = '82.50' },
    @{ Cell = 'E24'; Value = '  +1.84%  ' },
    @{ Cell = 'D25'; Value = '2.20' },
    @{ Cell = 'E25'; Value = '  -0.69%  ' },
    @{ Cell = 'D26'; Value = '11.98' },
    @{ Cell = 'E26'; Value = '  -3.60%  ' },
    @{ Cell = 'D27'; Value = '10.07' },
    @{ Cell = 'E27'; Value = '  -5.96%  ' },
    @{ Cell = 'E28'; Value = '  -0.02%  ' },
    @{ Cell = 'D29'; Value = '8.03' },
    @{ Cell = 'E29'; Value = '  +2.77%  ' },
    @{ Cell = 'D30'; Value = '2.39' },
    @{ Cell = 'E30'; Value = '  -2.25%  ' },
    @{ Cell = 'D31'; Value = '2.58' },
    @{ Cell = 'E31'; Value = '  -0.38%  ' },
    @{ Cell = 'D32'; Value = '0.0000101' },
    @{ Cell = 'E32'; Value = '  -1.05%  ' },
    @{ Cell = 'D33'; Value = '27.12' },
    @{ Cell = 'E33'; Value = '  +1.53%  ' },
    @{ Cell = 'D34'; Value = '0.110' },
    @{ Cell = 'E34'; Value = '  -1.46%  ' },
    @{ Cell = 'E35'; Value = '  -0.16%  ' },
    @{ Cell = 'D36'; Value = '0.978' },
    @{ Cell = 'E36'; Value = '  -0.70%  ' },
    @{ Cell = 'D37'; Value = '5.69' },
    @{ Cell = 'E37'; Value = '  -0.94%  ' },
    @{ Cell = 'D38'; Value = '49.02' },
    @{ Cell = 'E38'; Value = '  +0.03%  ' },
    @{ Cell = 'D39'; Value = '1.98' },
    @{ Cell = 'E39'; Value = '  -5.09%  ' },
    @{ Cell = 'D40'; Value = '43.98' },
    @{ Cell = 'E40'; Value = '  -4.24%  ' },
    @{ Cell = 'B41'; Value = 'dogwifhat' },
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif' },
    @{ Cell = 'D41'; Value = '2.85' },
    @{ Cell = 'E41'; Value = '  -1.78%  ' },
    @{ Cell = 'B42'; Value = 'TheGraph' },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt' },
    @{ Cell = 'D42'; Value = '0.299' },
    @{ Cell = 'E42'; Value = '  -0.44%  ' },
    @{ Cell = 'E43'; Value = '  -1.05%  ' },
    @{ Cell = 'D44'; Value = '8.39' },
    @{ Cell = 'E44'; Value = '  +0.12%  ' },
    @{ Cell = 'D45'; Value = '390.16' },
    @{ Cell = 'E45'; Value = '  +1.63%  ' },
    @{ Cell = 'E46'; Value = '  +0.62%  ' },
    @{ Cell = 'D47'; Value = '2.721.46' },
    @{ Cell = 'E47'; Value = '  -1.52%  ' },
    @{ Cell = 'D48'; Value = '131.14' },
    @{ Cell = 'E48'; Value = '  -2.54%  ' },
    @{ Cell = 'B50'; Value = 'InjectiveProtocol' },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' },
    @{ Cell = 'D50'; Value = '23.35' },
    @{ Cell = 'E50'; Value = '  +0.68%  ' },
    @{ Cell = 'B51'; Value = 'Stellar' },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm' },
    @{ Cell = 'D51'; Value = '0.106' },
    @{ Cell = 'E51'; Value = '  +1.01%  ' }
)

function Test-PlainNumber([string]$s) {
    # Matches simple unsigned integers/decimals such as "82.50" or
    # "0.0000101", but NOT multi-dot "thousands" strings like "64.788.15" -
    # those already round-trip through Excel as text untouched.
    return $s -match '^[0-9]+(\.[0-9]+)?$'
}

$textForced = @()

foreach ($u in $updates) {
    $cellRef = $u.Cell
    $value = $u.Value

    if ($cellRef.StartsWith("D") -and (Test-PlainNumber $value)) {
        # Force text so Excel doesn't auto-convert this into a Number cell.
        $ws.Range($cellRef).Value = "'" + $value
        $textForced += $cellRef
    } else {
        $ws.Range($cellRef).Value = $value
    }
}

# Strip the quote-prefix cell style picked up above so the cells end up with
# no special per-cell formatting, matching the rest of the sheet.
foreach ($ref in $textForced) {
    $ws.Range($ref).Style = "Normal"
}
